$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '64.133.64'
$c.Style = $s
$ws.Range("E2").Value = '  +0.60%  '
$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.145.03'
$c.Style = $s
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '591.69'
$c.Style = $s
$ws.Range("E5").Value = '  +0.37%  '
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '146.27'
$c.Style = $s
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  +0.09%  '
$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.136.61'
$c.Style = $s
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  -0.85%  '
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.Style = $s
$ws.Range("E10").Value = '  +0.36%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.Style = $s
$ws.Range("E11").Value = '  +2.07%  '
$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.459'
$c.Style = $s
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("E13").Value = '  -2.66%  '
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '37.21'
$c.Style = $s
$ws.Range("E14").Value = '  -0.11%  '
$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.664.50'
$c.Style = $s
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("E16").Value = '  -1.32%  '
$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.24'
$c.Style = $s
$ws.Range("E17").Value = '  +0.64%  '
$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '63.944.59'
$c.Style = $s
$ws.Range("E18").Value = '  +0.50%  '
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.143.11'
$c.Style = $s
$ws.Range("E19").Value = '  +0.52%  '
$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '466.42'
$c.Style = $s
$ws.Range("E20").Value = '  +0.06%  '
$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.36'
$c.Style = $s
$ws.Range("E22").Value = '  -0.42%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.49'
$c.Style = $s
$ws.Range("E23").Value = '  -0.90%  '
$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '13.03'
$c.Style = $s
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '81.13'
$c.Style = $s
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.32'
$c.Style = $s
$ws.Range("E26").Value = '  +5.67%  '
$ws.Range("E27").Value = '  +0.01%  '
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.65'
$c.Style = $s
$ws.Range("E28").Value = '  +7.52%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.26'
$c.Style = $s
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("E30").Value = '  +0.08%  '
$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.35'
$c.Style = $s
$ws.Range("E31").Value = '  +6.82%  '
$ws.Range("E32").Value = '  +0.17%  '
$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.48'
$c.Style = $s
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("E34").Value = '  +0.69%  '
$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0₃0831'
$c.Style = $s
$ws.Range("E35").Value = '  -5.85%  '
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  +0.58%  '
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = $s
$ws.Range("E38").Value = '  -2.68%  '
$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.24'
$c.Style = $s
$ws.Range("E39").Value = '  -5.81%  '
$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '459.33'
$c.Style = $s
$ws.Range("E40").Value = '  +1.18%  '
$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '51.30'
$c.Style = $s
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("E42").Value = '  +5.31%  '
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("E44").Value = '  -0.16%  '
$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.922.53'
$c.Style = $s
$ws.Range("E45").Value = '  +0.63%  '
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '39.51'
$c.Style = $s
$ws.Range("E46").Value = '  +9.48%  '
$ws.Range("E47").Value = '  -2.50%  '
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '130.63'
$c.Style = $s
$ws.Range("E48").Value = '  +2.43%  '
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = $s
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("E51").Value = '  -1.01%  '
